# Add a "Value" column (C) with numbers, a SUBTOTAL row below the data,
# and extend the autofilter / filter-database range to cover it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 starts out hidden (filtered out). Temporarily unhide it while we
# write into it so the engine doesn't stamp a stray custom row height,
# then restore the hidden state.
$ws.Rows.Item(3).Hidden = $false

# New "Value" column header + data
$ws.Range("C1").Value = "Value"
$ws.Range("C2").Value = 5
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 15

$ws.Rows.Item(3).Hidden = $true

# Subtotal formula below the data (sums the visible rows of column C)
$ws.Range("C6").Formula = "=SUBTOTAL(9,C2:C4)"

# Re-apply the autofilter over the extended A1:C4 range, keeping the
# existing criteria (values 1 and 3) on the first column
$ws.AutoFilterMode = $false
$ws.Range("A1:C4").AutoFilter(1, @("1","3"), 7)

# Update the hidden _xlnm._FilterDatabase defined name to match the new range
$wb.Names.Item(1).RefersTo = '=Filtered!$A$1:$C$4'

# Move the selection to match the edited workbook
$ws.Range("D6").Select()
